$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 7708
$ws.Range("K2").Value = 2353
$ws.Range("J3").Value = 8077
$ws.Range("K3").Value = 2267
$ws.Range("K4").Value = 480
$ws.Range("K5").Value = 149
$ws.Range("K6").Value = 2840
$ws.Range("K7").Value = 8089

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 160
$ws.Range("K3").Value = 158
$ws.Range("K4").Value = 30
$ws.Range("K6").Value = 181
$ws.Range("K7").Value = 541

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 56
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 44
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 90
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 260

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 55
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 239
$ws.Range("K8").Value = 541
$ws.Range("K10").Value = 46
$ws.Range("K11").Value = 176
$ws.Range("K15").Value = 81
$ws.Range("K16").Value = 19
$ws.Range("K19").Value = 233
$ws.Range("K20").Value = 178
$ws.Range("K23").Value = 73
$ws.Range("K24").Value = 24
$ws.Range("K27").Value = 88
$ws.Range("K29").Value = 412
$ws.Range("K31").Value = 92
$ws.Range("K32").Value = 14
$ws.Range("K33").Value = 315
$ws.Range("K36").Value = 97
$ws.Range("K37").Value = 260
$ws.Range("K42").Value = 280
$ws.Range("K43").Value = 74
$ws.Range("K44").Value = 77
$ws.Range("K48").Value = 99
$ws.Range("K52").Value = 220
$ws.Range("K53").Value = 115
$ws.Range("K54").Value = 154
$ws.Range("K57").Value = 24
$ws.Range("K64").Value = 54
$ws.Range("K65").Value = 194
$ws.Range("K67").Value = 313
$ws.Range("K69").Value = 21
$ws.Range("K71").Value = 22
$ws.Range("K72").Value = 37
$ws.Range("K74").Value = 7
$ws.Range("K76").Value = 116
$ws.Range("K79").Value = 213
$ws.Range("K83").Value = 177
$ws.Range("K85").Value = 394
$ws.Range("K86").Value = 56
$ws.Range("K90").Value = 70
$ws.Range("K91").Value = 75
$ws.Range("K94").Value = 97
$ws.Range("K95").Value = 129
$ws.Range("K96").Value = 110
$ws.Range("K99").Value = 145
$ws.Range("K100").Value = 12
$ws.Range("K101").Value = 8089

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 98
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 313

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 136
$ws.Range("K4").Value = 23
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 412

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 20
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 71
$ws.Range("K3").Value = 83
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 69
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 213

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 60
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 38
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 259
$ws.Range("K2").Value = 78
$ws.Range("J3").Value = 252
$ws.Range("K3").Value = 77
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 54
$ws.Range("K3").Value = 45
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 134
$ws.Range("K7").Value = 394

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 60
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 7

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Range("K4").Value = 1
$ws.Range("K6").Value = 4
